$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.310.77"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.609.51"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.20"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.48"
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.833.14"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.596.35"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.283.54"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.25"
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.07"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.35"
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.03"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.89"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.42"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.24"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("E30").Value = "  +5.29%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("E32").Value = "  +2.80%  "
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.163.27"
$ws.Range("E36").Value = "  +3.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0167"
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.790"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("E42").Value = "  +4.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.783"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.744.32"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.22"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("E46").Value = "  +11.84%  "
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.92"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  -0.17%  "
